$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert a new row above row 5 (most recent transaction goes on top,
# pushing the existing rows 5 and 6 down to 6 and 7 respectively).
$ws.Rows.Item(5).Insert()

# The freshly inserted row inherits the bold/filled header formatting
# (and the header's full A:AB column span) from the row above it;
# reset it back to the default "Normal" style before re-applying the
# date number format, so we land on the exact same style that was
# already used for the date column elsewhere, then drop the spurious
# empty cells beyond column J that the insert pulled in from the header.
$ws.Range("A5:AB5").Style = "Normal"
$ws.Range("K5:AB5").Clear()
$ws.Cells.Item(5, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"

# Populate the new transaction row.
$ws.Cells.Item(5, 1).Value = 46066
$ws.Cells.Item(5, 2).Value = "NSE"
$ws.Cells.Item(5, 3).Value = "Buy"
$ws.Cells.Item(5, 4).Value = 30
$ws.Cells.Item(5, 5).Value = 140.47
$ws.Cells.Item(5, 6).Value = 4243.72
$ws.Cells.Item(5, 7).Value = "CN#252611910666"
$ws.Cells.Item(5, 8).Value = 4.2567
$ws.Cells.Item(5, 9).Value = 25.3647
$ws.Cells.Item(5, 10).Formula = "=Index!`$C`$2"

Write-Output "Row inserted and populated"
